$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-28 Saturday" "2024-09-29 Sunday"

Replace-Text "867×4=" "987×3="
Replace-Text "817×7=" "312×8="
Replace-Text "688×3=" "689×6="
Replace-Text "937×7=" "801×5="
Replace-Text "755×2=" "832×4="

Replace-Text "286×7=" "970×3="
Replace-Text "160×4=" "716×8="
Replace-Text "376×6=" "870×4="
Replace-Text "178×9=" "253×7="
Replace-Text "441×7=" "856×3="

Replace-Text "445×3=" "907×9="
Replace-Text "535×2=" "938×8="
Replace-Text "738×4=" "209×8="
Replace-Text "121×8=" "558×6="
Replace-Text "734×2=" "936×6="

Replace-Text "709×9=" "817×6="
Replace-Text "297×3=" "679×6="
Replace-Text "691×6=" "743×5="
Replace-Text "235×8=" "501×9="
Replace-Text "622×3=" "989×7="

Replace-Text "417×9=" "304×7="
Replace-Text "976×4=" "807×9="
Replace-Text "500×6=" "401×9="
Replace-Text "357×2=" "496×7="
Replace-Text "449×7=" "164×5="
